# "add last plast of pochistone 019"
# - Mark plates 016/017/018 as DONE (rename tabs).
# - Fill in the last (third) row of 12 wells for plate 019's layout
#   (rows 2 & 3 get the remaining wells finished, row 4 is brand new).
# - Make "PocHistone RLFP 019" the active/selected sheet with the cursor
#   parked at H14 (where the next entry would go).

$wb = $excel.ActiveWorkbook

# --- Rename the now-finished plates -----------------------------------
$wb.Worksheets.Item("PocHistone RLFP 016").Name = "DONE PocHistone RLFP 016"
$wb.Worksheets.Item("PocHistone RLFP 017").Name = "DONE PocHistone RLFP 017"
$wb.Worksheets.Item("PocHistone RLFP 018").Name = "DONE PocHistone RLFP 018"

# --- Update plate 019's well-assignment grid ---------------------------
$ws = $wb.Worksheets.Item("PocHistone RLFP 019")

# Row "A" (row 2) - finish filling out the remaining wells
$ws.Range("C2").Value = "P13 E11 480"
$ws.Range("D2").Value = "P13 E12 487"
$ws.Range("E2").Value = "P13 D11 613"
$ws.Range("F2").Value = "P5 G8 760"
$ws.Range("G2").Value = "P37 G9 802"
$ws.Range("H2").Value = "P6 F6 1011"
$ws.Range("I2").Value = "P5 H1 1098"
$ws.Range("J2").Value = "P38 D11 3940"
$ws.Range("K2").Value = "P37 A3 1160"
$ws.Range("L2").Value = "P19 A4 1321"
$ws.Range("M2").Value = "P19 C1 1418"

# Row "B" (row 3) - finish filling out the remaining wells
$ws.Range("B3").Value = "+"
$ws.Range("C3").Value = "P16 H4 1444"
$ws.Range("D3").Value = "P35 A11 1456"
$ws.Range("E3").Value = "P34 C8 1514"
$ws.Range("F3").Value = "P34 C2 1582"
$ws.Range("G3").Value = "P4 C8 1676"
$ws.Range("H3").Value = "P33 A3 1836"
$ws.Range("I3").Value = "P32 F3 1840"
$ws.Range("J3").Value = "P29 H8 2437"
$ws.Range("K3").Value = "P30 C6 2555"
$ws.Range("L3").Value = "P29 C10 2604"
$ws.Range("M3").Value = "P30 E10 2676"

# Row "C" (row 4) - brand new row of samples
$ws.Range("A4").Value = "C"
$ws.Range("B4").Value = "P26 D6 2867"
$ws.Range("C4").Value = "P21 F6 3102"
$ws.Range("D4").Value = "P20 F11 3242"
$ws.Range("E4").Value = "P28 A9 3500"
$ws.Range("F4").Value = "P24 A9 3507"
$ws.Range("G4").Value = "P28 B12 3527"
$ws.Range("H4").Value = "P14 A11 3729"
$ws.Range("I4").Value = "P38 D4 3910"
$ws.Range("J4").Value = "P12 G4 3919"
$ws.Range("K4").Value = "P38 D7 3928"
$ws.Range("L4").Value = "P38 D9 3830"
$ws.Range("M4").Value = "P38 C9 3936"

# --- Leave the workbook focused on plate 019 ---------------------------
$ws.Activate()
$ws.Range("H14").Select()
